# New version of the 2nd lead experiment: rows 141-199 in the "lead_v2"
# sheet's column B (time) were background/no-absorber runs, so relabel
# the numeric placeholder 0 with the text marker "bg".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lead_v2")

# Replace the old numeric 0 placeholders in B141:B199 with the text "bg".
$ws.Range("B141:B199").Value = "bg"

# Make sure lead_v2 is the active sheet, then move the selection to the
# range that was just edited (matches the saved cursor position).
$ws.Activate() | Out-Null
$ws.Range("B141:B199").Select() | Out-Null

# Reposition the Excel window the way it was when the file was saved.
$excel.ActiveWindow.Left = 28005
$excel.ActiveWindow.Top = 3885
